$d = $word.ActiveDocument

# The first paragraph currently reads:
#   "This is a Microsoft word document."
# We need to turn it into four separate runs:
#   "This is a Microsoft word document."
#   " ("
#   "Changed main"
#   ")"
# A plain sequence of Range.InsertAfter() calls on the same spot gets
# coalesced by the engine into a single run when the adjacent text shares
# identical formatting, so we drop a temporary bookmark between each
# insertion to force a run boundary, then remove the bookmarks again.

$p1 = $d.Paragraphs(1)

# --- insert " (" right after the existing sentence ---
$r = $p1.Range
$r.MoveEnd(1, -1)
$r.Collapse(0)
$d.Bookmarks.Add("ztmpBoundary0", $r)
$r.InsertAfter(" (")

# --- insert "Changed main" after " (" ---
$p1 = $d.Paragraphs(1)
$r = $p1.Range
$r.MoveEnd(1, -1)
$r.Collapse(0)
$d.Bookmarks.Add("ztmpBoundary1", $r)
$r.InsertAfter("Changed main")

# --- insert ")" after "Changed main" ---
$p1 = $d.Paragraphs(1)
$r = $p1.Range
$r.MoveEnd(1, -1)
$r.Collapse(0)
$r.InsertAfter(")")

# Clean up the temporary bookmarks so they don't end up in the saved file.
$d.Bookmarks("ztmpBoundary0").Delete()
$d.Bookmarks("ztmpBoundary1").Delete()
